# Economic Dashboard update - 2025-12-10
# Applies the refreshed release dates / rolling-window values described in
# the commit "Update dashboards - 2025-12-10".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Style-only changes: the "N" (and one "C") release-date cells toggle
# between the plain date format (style of C5 / "donor47") and the
# yellow "just refreshed" highlight (style of C7 / "donor48").
# Using Copy + PasteSpecial(xlPasteFormats) re-uses the workbook's
# existing style entries instead of minting new ones.
# ---------------------------------------------------------------------
$xlPasteFormats = -4122

function Copy-Style($fromAddr, $toAddr) {
    $ws.Range($fromAddr).Copy()
    $ws.Range($toAddr).PasteSpecial($xlPasteFormats)
}

# N5: highlighted -> normal (value unchanged)
Copy-Style "C5" "N5"

# N10/N11/N12: normal -> highlighted (JOLTS block refreshed)
Copy-Style "C7" "N10"
Copy-Style "C7" "N11"
Copy-Style "C7" "N12"

# N31/N32: normal -> highlighted (ECI wages refreshed)
Copy-Style "C7" "N31"
Copy-Style "C7" "N32"

# C32/C33/C34: highlighted -> normal
Copy-Style "C5" "C32"
Copy-Style "C5" "C33"
Copy-Style "C5" "C34"

# N41/N42/N43/N44: highlighted -> normal
Copy-Style "C5" "N41"
Copy-Style "C5" "N42"
Copy-Style "C5" "N43"
Copy-Style "C5" "N44"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Row 5 - ADP Total NonFarm Private: release date refreshed (no new obs)
# ---------------------------------------------------------------------
$ws.Range("N5").Value = 45962

# ---------------------------------------------------------------------
# Row 10 - JOLTS Openings Rate: new observation, window rolls forward
# ---------------------------------------------------------------------
$ws.Range("N10").Value = 45931
$ws.Range("Q10").Value = 4.6
$ws.Range("R10").Value = 4.6
$ws.Range("S10").Value = 4.3
$ws.Range("T10").Value = 4.3

# ---------------------------------------------------------------------
# Row 11 - JOLTS Hires Rate: new observation, window rolls forward
# ---------------------------------------------------------------------
$ws.Range("N11").Value = 45931
$ws.Range("R11").Value = 3.4
$ws.Range("S11").Value = 3.2
$ws.Range("T11").Value = 3.3
$ws.Range("U11").Value = 3.3

# ---------------------------------------------------------------------
# Row 12 - JOLTS Separations Rate: new observation, window rolls forward
# ---------------------------------------------------------------------
$ws.Range("N12").Value = 45931
$ws.Range("S12").Value = 3.2

# ---------------------------------------------------------------------
# Row 29 - 5yr,5yr Forward: new daily observation
# ---------------------------------------------------------------------
$ws.Range("N29").Value = 46000
$ws.Range("R29").Value = 2.2
$ws.Range("S29").Value = 0
$ws.Range("T29").Value = 0

# ---------------------------------------------------------------------
# Row 30 - 10yr TIPS: new daily observation
# ---------------------------------------------------------------------
$ws.Range("N30").Value = 46000
$ws.Range("S30").Value = 0
$ws.Range("T30").Value = 0
$ws.Range("U30").Value = 2.26

# ---------------------------------------------------------------------
# Row 31 - ECI Wages Q/Q: new quarterly observation, window rolls forward
# ---------------------------------------------------------------------
$ws.Range("N31").Value = 45839
$ws.Range("Q31").Value = 0.007996957929548465
$ws.Range("R31").Value = 0.01027939464493599
$ws.Range("S31").Value = 0.007624633431085215
$ws.Range("T31").Value = 0.009473060982829962
$ws.Range("U31").Value = 0.007756563245823411

# ---------------------------------------------------------------------
# Row 32 - ECI Wages Y/Y: new quarterly observation, window rolls forward
# ---------------------------------------------------------------------
$ws.Range("N32").Value = 45839
$ws.Range("Q32").Value = 0.03584369449378332
$ws.Range("R32").Value = 0.03559665871121723
$ws.Range("S32").Value = 0.03369434416365838
$ws.Range("T32").Value = 0.03710462287104619
$ws.Range("U32").Value = 0.03746928746928743

# ---------------------------------------------------------------------
# Row 47 - Fed Funds Rate: release date refreshed (no new obs)
# ---------------------------------------------------------------------
$ws.Range("N47").Value = 45999

# ---------------------------------------------------------------------
# Row 48 - 2y UST: new daily observation, window rolls forward
# ---------------------------------------------------------------------
$ws.Range("N48").Value = 45999
$ws.Range("Q48").Value = 3.57
$ws.Range("R48").Value = 0
$ws.Range("S48").Value = 0
$ws.Range("T48").Value = 3.56
$ws.Range("U48").Value = 3.52

# ---------------------------------------------------------------------
# Row 49 - 5y UST: new daily observation, window rolls forward
# ---------------------------------------------------------------------
$ws.Range("N49").Value = 45999
$ws.Range("Q49").Value = 3.75
$ws.Range("R49").Value = 0
$ws.Range("S49").Value = 0
$ws.Range("T49").Value = 3.72
$ws.Range("U49").Value = 3.68

# ---------------------------------------------------------------------
# Row 50 - 10y UST: new daily observation, window rolls forward
# ---------------------------------------------------------------------
$ws.Range("N50").Value = 45999
$ws.Range("Q50").Value = 4.17
$ws.Range("R50").Value = 0
$ws.Range("S50").Value = 0
$ws.Range("T50").Value = 4.14
$ws.Range("U50").Value = 4.11

# ---------------------------------------------------------------------
# Row 52 - BAA: new daily observation, window rolls forward
# ---------------------------------------------------------------------
$ws.Range("N52").Value = 45999
$ws.Range("Q52").Value = 5.9
$ws.Range("R52").Value = 0
$ws.Range("S52").Value = 0
$ws.Range("T52").Value = 5.88
